$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updated odds
$ws.Cells.Item(2, 8).Value = 3.35
$ws.Cells.Item(2, 9).Value = 4.15
$ws.Cells.Item(2, 10).Value = 2.4
$ws.Cells.Item(2, 11).Value = 2.12
$ws.Cells.Item(2, 12).Value = 4.55
$ws.Cells.Item(2, 13).Value = 1.07
$ws.Cells.Item(2, 14).Value = 6.8
$ws.Cells.Item(2, 15).Value = 1.34
$ws.Cells.Item(2, 16).Value = 3
$ws.Cells.Item(2, 17).Value = 2
$ws.Cells.Item(2, 18).Value = 1.72
$ws.Cells.Item(2, 19).Value = 1.42
$ws.Cells.Item(2, 20).Value = 2.67
$ws.Cells.Item(2, 21).Value = 1.87
$ws.Cells.Item(2, 22).Value = 1.85
$ws.Cells.Item(2, 23).Value = 6.6
$ws.Cells.Item(2, 24).Value = 8.5
$ws.Cells.Item(2, 25).Value = 8.25
$ws.Cells.Item(2, 26).Value = 15.5
$ws.Cells.Item(2, 27).Value = 15
$ws.Cells.Item(2, 28).Value = 28
$ws.Cells.Item(2, 29).Value = 6.8
$ws.Cells.Item(2, 30).Value = 6.4
$ws.Cells.Item(2, 31).Value = 15.5
$ws.Cells.Item(2, 32).Value = 80
$ws.Cells.Item(2, 33).Value = 700
$ws.Cells.Item(2, 34).Value = 10.5
$ws.Cells.Item(2, 35).Value = 22
$ws.Cells.Item(2, 39).Value = 50
$ws.Cells.Item(2, 40).Value = 3.7
$ws.Cells.Item(2, 41).Value = 9
$ws.Cells.Item(2, 42).Value = 18
$ws.Cells.Item(2, 43).Value = 32
$ws.Cells.Item(2, 44).Value = 65
$ws.Cells.Item(2, 45).Value = 250
$ws.Cells.Item(2, 46).Value = 2.67
$ws.Cells.Item(2, 47).Value = 7.3
$ws.Cells.Item(2, 48).Value = 70
$ws.Cells.Item(2, 50).Value = 24
$ws.Cells.Item(2, 51).Value = 30
$ws.Cells.Item(2, 54).Value = 400

# Row 3 updated odds
$ws.Cells.Item(3, 7).Value = 2.45
$ws.Cells.Item(3, 9).Value = 2.88
$ws.Cells.Item(3, 10).Value = 3.2
$ws.Cells.Item(3, 11).Value = 2.05
$ws.Cells.Item(3, 17).Value = 2.08
$ws.Cells.Item(3, 18).Value = 1.73
$ws.Cells.Item(3, 23).Value = 7.5
$ws.Cells.Item(3, 25).Value = 10
$ws.Cells.Item(3, 29).Value = 9
$ws.Cells.Item(3, 33).Value = 301
$ws.Cells.Item(3, 34).Value = 8.5
$ws.Cells.Item(3, 35).Value = 13
$ws.Cells.Item(3, 40).Value = 4.5
$ws.Cells.Item(3, 41).Value = 15
$ws.Cells.Item(3, 43).Value = 51
$ws.Cells.Item(3, 45).Value = 201

# Row 6 updated odds
$ws.Cells.Item(6, 17).Value = 1.57
$ws.Cells.Item(6, 18).Value = 2.35

# Row 7 - new match (Almeria vs Cordoba)
$ws.Cells.Item(7, 1).Value = "MyS7z3gR"
$ws.Cells.Item(7, 2).Value = "26/11/2024"
$ws.Cells.Item(7, 3).Value = "16:30"
$ws.Cells.Item(7, 4).Value = "SPAIN - LALIGA2"
$ws.Cells.Item(7, 5).Value = "Almeria"
$ws.Cells.Item(7, 6).Value = "Cordoba"
$ws.Cells.Item(7, 7).Value = 1.65
$ws.Cells.Item(7, 8).Value = 4.1
$ws.Cells.Item(7, 9).Value = 4.75
$ws.Cells.Item(7, 10).Value = 2.2
$ws.Cells.Item(7, 11).Value = 2.4
$ws.Cells.Item(7, 12).Value = 4.75
$ws.Cells.Item(7, 13).Value = 1.03
$ws.Cells.Item(7, 14).Value = 15
$ws.Cells.Item(7, 15).Value = 1.18
$ws.Cells.Item(7, 16).Value = 4.5
$ws.Cells.Item(7, 17).Value = 1.62
$ws.Cells.Item(7, 18).Value = 2.25
$ws.Cells.Item(7, 19).Value = 1.3
$ws.Cells.Item(7, 20).Value = 3.4
$ws.Cells.Item(7, 21).Value = 1.67
$ws.Cells.Item(7, 22).Value = 2.1
$ws.Cells.Item(7, 23).Value = 9
$ws.Cells.Item(7, 24).Value = 9
$ws.Cells.Item(7, 25).Value = 8.5
$ws.Cells.Item(7, 26).Value = 13
$ws.Cells.Item(7, 27).Value = 12
$ws.Cells.Item(7, 28).Value = 21
$ws.Cells.Item(7, 29).Value = 15
$ws.Cells.Item(7, 30).Value = 8
$ws.Cells.Item(7, 31).Value = 13
$ws.Cells.Item(7, 32).Value = 41
$ws.Cells.Item(7, 33).Value = 151
$ws.Cells.Item(7, 34).Value = 17
$ws.Cells.Item(7, 35).Value = 26
$ws.Cells.Item(7, 36).Value = 15
$ws.Cells.Item(7, 37).Value = 51
$ws.Cells.Item(7, 38).Value = 34
$ws.Cells.Item(7, 39).Value = 34
$ws.Cells.Item(7, 40).Value = 3.75
$ws.Cells.Item(7, 41).Value = 8.5
$ws.Cells.Item(7, 42).Value = 17
$ws.Cells.Item(7, 43).Value = 23
$ws.Cells.Item(7, 44).Value = 41
$ws.Cells.Item(7, 45).Value = 101
$ws.Cells.Item(7, 46).Value = 3.4
$ws.Cells.Item(7, 47).Value = 7.5
$ws.Cells.Item(7, 48).Value = 51
$ws.Cells.Item(7, 49).Value = 6.5
$ws.Cells.Item(7, 50).Value = 23
$ws.Cells.Item(7, 51).Value = 29
$ws.Cells.Item(7, 52).Value = 81
$ws.Cells.Item(7, 53).Value = 81
$ws.Cells.Item(7, 54).Value = 151
$ws.Cells.Item(7, 55).Value = 81
$ws.Cells.Item(7, 56).Value = 81

# Row 8 - new match (Rampla Juniors vs Liverpool M.)
$ws.Cells.Item(8, 1).Value = "tznSwPz8"
$ws.Cells.Item(8, 2).Value = "26/11/2024"
$ws.Cells.Item(8, 3).Value = "16:30"
$ws.Cells.Item(8, 4).Value = "URUGUAY - PRIMERA DIVISION"
$ws.Cells.Item(8, 5).Value = "Rampla Juniors"
$ws.Cells.Item(8, 6).Value = "Liverpool M."
$ws.Cells.Item(8, 7).Value = 2.63
$ws.Cells.Item(8, 8).Value = 3.1
$ws.Cells.Item(8, 9).Value = 2.75
$ws.Cells.Item(8, 10).Value = 3.4
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 3.5
$ws.Cells.Item(8, 13).Value = 1.07
$ws.Cells.Item(8, 14).Value = 9
$ws.Cells.Item(8, 15).Value = 1.36
$ws.Cells.Item(8, 16).Value = 3
$ws.Cells.Item(8, 17).Value = 2.2
$ws.Cells.Item(8, 18).Value = 1.65
$ws.Cells.Item(8, 19).Value = 1.5
$ws.Cells.Item(8, 20).Value = 2.5
$ws.Cells.Item(8, 21).Value = 1.83
$ws.Cells.Item(8, 22).Value = 1.83
$ws.Cells.Item(8, 23).Value = 7.5
$ws.Cells.Item(8, 24).Value = 12
$ws.Cells.Item(8, 25).Value = 11
$ws.Cells.Item(8, 26).Value = 26
$ws.Cells.Item(8, 27).Value = 23
$ws.Cells.Item(8, 28).Value = 34
$ws.Cells.Item(8, 29).Value = 8
$ws.Cells.Item(8, 30).Value = 6
$ws.Cells.Item(8, 31).Value = 15
$ws.Cells.Item(8, 32).Value = 51
$ws.Cells.Item(8, 33).Value = 351
$ws.Cells.Item(8, 34).Value = 8
$ws.Cells.Item(8, 35).Value = 13
$ws.Cells.Item(8, 36).Value = 11
$ws.Cells.Item(8, 37).Value = 29
$ws.Cells.Item(8, 38).Value = 23
$ws.Cells.Item(8, 39).Value = 34
$ws.Cells.Item(8, 40).Value = 4.5
$ws.Cells.Item(8, 41).Value = 15
$ws.Cells.Item(8, 42).Value = 26
$ws.Cells.Item(8, 43).Value = 51
$ws.Cells.Item(8, 44).Value = 81
$ws.Cells.Item(8, 45).Value = 201
$ws.Cells.Item(8, 46).Value = 2.5
$ws.Cells.Item(8, 47).Value = 8.5
$ws.Cells.Item(8, 48).Value = 67
$ws.Cells.Item(8, 49).Value = 4.75
$ws.Cells.Item(8, 50).Value = 17
$ws.Cells.Item(8, 51).Value = 29
$ws.Cells.Item(8, 52).Value = 51
$ws.Cells.Item(8, 53).Value = 81
$ws.Cells.Item(8, 54).Value = 201
$ws.Cells.Item(8, 55).Value = 51
$ws.Cells.Item(8, 56).Value = 51

